$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new sheets ("Putz 205" and "Putz 206 APS") at the end, right
#    after the current last sheet ("Putz 204").
# ---------------------------------------------------------------------------
$lastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)

$sheet205 = $wb.Worksheets.Add($null, $lastExisting)
$sheet205.Name = "Putz 205"

$sheet206 = $wb.Worksheets.Add($null, $sheet205)
$sheet206.Name = "Putz 206 APS"

# ---------------------------------------------------------------------------
# Helper-ish block: build "Putz 205" content
# ---------------------------------------------------------------------------
$ws = $sheet205

$ws.Columns.Item(1).ColumnWidth = 6.29
$ws.Columns.Item(2).ColumnWidth = 12.71
$ws.Columns.Item(4).ColumnWidth = 21.86

# Header row
$ws.Range("A1").Value = "Sr. No"
$ws.Range("B1").Value = "PO NO."
$ws.Range("C1").Value = "Part No."
$ws.Range("D1").Value = "Item Description"
$ws.Range("E1").Value = "Qty"
$ws.Range("F1").Value = "Price"
$ws.Range("G1").Value = "Total"
$ws.Rows.Item(1).RowHeight = 21.6

$hdr = $ws.Range("A1:G1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# Row 2 - line item 1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "151357107 (25-07-2023)"
$ws.Range("C2").Value = 615181
$ws.Range("D2").Value = "INDUSTRIAL BELL_225mm"
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 2400
$ws.Range("G2").Formula = "=E2*F2"
$ws.Rows.Item(2).RowHeight = 37.8

# Row 3 - line item 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "151362092 (08-08-2023)"
$ws.Range("C3").Value = 616026
$ws.Range("D3").Value = "Screw Conveyor Bellow_Dia 224x500mm"
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 5614
$ws.Range("G3").Formula = "=E3*F3"
$ws.Rows.Item(3).RowHeight = 43.8

$data = $ws.Range("A2:G3")
$data.Borders.LineStyle = 1
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108
$data.WrapText = $true

# Row 4 - Total
$ws.Range("A4").Value = "Total"
$ws.Range("A4:F4").Merge()
$ws.Range("G4").Formula = "=SUM(G2:G3)"

# Row 5 - CGST 9%
$ws.Range("A5").Value = "CGST 9%"
$ws.Range("A5:F5").Merge()
$ws.Range("G5").Formula = "=G4*9%"

# Row 6 - SGST 9%
$ws.Range("A6").Value = "SGST 9%"
$ws.Range("A6:F6").Merge()
$ws.Range("G6").Formula = "=G4*9%"

# Row 7 - Grand Total
$ws.Range("A7").Value = "Grand Total"
$ws.Range("A7:F7").Merge()
$ws.Range("G7").Formula = "=SUM(G4:G6)"

$totals = $ws.Range("A4:G7")
$totals.Font.Bold = $true
$totals.Borders.LineStyle = 1
$totals.HorizontalAlignment = -4108
$totals.VerticalAlignment = -4108
$totals.WrapText = $true

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# Helper-ish block: build "Putz 206 APS" content
# ---------------------------------------------------------------------------
$ws2 = $sheet206

$ws2.Columns.Item(1).ColumnWidth = 7.51
$ws2.Columns.Item(2).ColumnWidth = 13.38
$ws2.Columns.Item(4).ColumnWidth = 21.75

# Header row
$ws2.Range("A1").Value = "Sr. No"
$ws2.Range("B1").Value = "PO NO."
$ws2.Range("C1").Value = "Part No."
$ws2.Range("D1").Value = "Item Description"
$ws2.Range("E1").Value = "Qty"
$ws2.Range("F1").Value = "Price"
$ws2.Range("G1").Value = "Total"
$ws2.Rows.Item(1).RowHeight = 29.4

$hdr2 = $ws2.Range("A1:G1")
$hdr2.Font.Bold = $true
$hdr2.Borders.LineStyle = 1
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4108
$hdr2.WrapText = $true

# Row 2 - line item 1
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "APS - 151362001 (08-08-2023)"
$ws2.Range("C2").Value = 616026
$ws2.Range("D2").Value = "Screw Conveyor Bellow_Dia 224x500mm"
$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 5614
$ws2.Range("G2").Formula = "=E2*F2"
$ws2.Rows.Item(2).RowHeight = 42.6

$data2 = $ws2.Range("A2:G2")
$data2.Borders.LineStyle = 1
$data2.HorizontalAlignment = -4108
$data2.VerticalAlignment = -4108
$data2.WrapText = $true

# Row 3 - Total
$ws2.Range("A3").Value = "Total"
$ws2.Range("A3:F3").Merge()
$ws2.Range("G3").Formula = "=SUM(G2)"

# Row 4 - CGST 9%
$ws2.Range("A4").Value = "CGST 9%"
$ws2.Range("A4:F4").Merge()
$ws2.Range("G4").Formula = "=G3*9%"

# Row 5 - SGST 9%
$ws2.Range("A5").Value = "SGST 9%"
$ws2.Range("A5:F5").Merge()
$ws2.Range("G5").Formula = "=G3*9%"

# Row 6 - Grand Total
$ws2.Range("A6").Value = "Grand Total"
$ws2.Range("A6:F6").Merge()
$ws2.Range("G6").Formula = "=SUM(G3:G5)"

$totals2 = $ws2.Range("A3:G6")
$totals2.Font.Bold = $true
$totals2.Borders.LineStyle = 1
$totals2.HorizontalAlignment = -4108
$totals2.VerticalAlignment = -4108
$totals2.WrapText = $true

$ws2.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. "Putz 204" is no longer the tab-selected sheet; select its full used
#    range instead of the old ad-hoc C16 cell.
# ---------------------------------------------------------------------------
$sheet204 = $wb.Worksheets.Item("Putz 204")
$sheet204.Range("A1:G8").Select()

# ---------------------------------------------------------------------------
# 3. "Putz 206 APS" becomes the active / tab-selected sheet, with its whole
#    used range selected.
# ---------------------------------------------------------------------------
$sheet206.Activate()
$sheet206.Range("A1:G6").Select()

# ---------------------------------------------------------------------------
# 4. Scroll / active-tab bookkeeping on the workbook window.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollWorkbookTabs(0)
for ($i = 1; $i -le 22; $i++) {
    $excel.ActiveWindow.ScrollWorkbookTabs(1)
}
